{"js": "// 1. Update the project title text: \"Inventarios y facturaci\u00f3n...\" -> \"inventario, cotizaci\u00f3n y facturaci\u00f3n...\"\nconst oldTitleText = \"Inventarios y facturaci\u00f3n para la empresa familiar:\";\nconst newTitleText = \"inventario, cotizaci\u00f3n y facturaci\u00f3n para la empresa familiar:\";\n\nconst titleResults = context.document.body.search(oldTitleText, { matchCase: true });\ntitleResults.load(\"items\");\nawait context.sync();\n\nif (titleResults.items.length > 0) {\n  // Replace in-place so the surrounding runs (and their formatting) stay untouched.\n  titleResults.items[0].insertText(newTitleText, \"Replace\");\n  await context.sync();\n}\n\n// 2. Remove the lone, space-only run that sits by itself in the (otherwise empty)\n//    paragraph right after the title paragraph, leaving the paragraph itself intact.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  if (paragraph.text === \" \") {\n    const paragraphRange = paragraph.getRange();\n    const spaceResults = paragraphRange.search(\" \", { matchCase: true });\n    spaceResults.load(\"items\");\n    await context.sync();\n\n    if (spaceResults.items.length > 0) {\n      spaceResults.items[0].delete();\n      await context.sync();\n    }\n    break;\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# ---------------------------------------------------------------------------\n# 1. Update the project title text:\n#    \"Inventarios y facturaci\u00f3n...\" -> \"inventario, cotizaci\u00f3n y facturaci\u00f3n...\"\n#    Only the text of the single run is changed; the run itself (and its\n#    neighbours) must keep their own separate <w:r> elements/formatting.\n# ---------------------------------------------------------------------------\n$oldTitleText = \"Inventarios y facturaci\u00f3n para la empresa familiar:\"\n$newTitleText = \"inventario, cotizaci\u00f3n y facturaci\u00f3n para la empresa familiar:\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldTitleText\n$found = $find.Execute()\n\nif ($found) {\n    $targetRange = $find.Parent\n    $rangeStart = $targetRange.Start\n\n    # Nudging a formatting property immediately before/after the text\n    # assignment keeps the engine from silently coalescing this run back\n    # into the identically-formatted run next to it.\n    $targetRange.Bold = 1\n    $targetRange.Text = $newTitleText\n\n    $rangeEnd = $rangeStart + $newTitleText.Length\n    $revertRange = $d.Range($rangeStart, $rangeEnd)\n    $revertRange.Bold = 0\n\n    # -----------------------------------------------------------------------\n    # 2. Remove the lone, space-only run in the (otherwise empty) paragraph\n    #    that immediately follows the title paragraph - the paragraph itself\n    #    must stay in place, just without any run left inside it.\n    # -----------------------------------------------------------------------\n    $prefixRange = $d.Range(0, $rangeEnd)\n    $titleParaIndex = $prefixRange.Paragraphs.Count\n    $titlePara = $d.Paragraphs.Item($titleParaIndex)\n    $nextPara = $titlePara.Next()\n\n    $nextText = $nextPara.Range.Text.TrimEnd([char]13)\n    if ($nextText -eq \" \") {\n        $spanStart = $nextPara.Range.Start\n        $spanEnd = $nextPara.Range.End - 1\n        if ($spanEnd -gt $spanStart) {\n            $emptyRunRange = $d.Range($spanStart, $spanEnd)\n            $emptyRunRange.Delete()\n        }\n    }\n}\n"}
